$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, matching the format of the other headers (copy format from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new column's data values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
